# Edit: insert a new data row at row 95 (shifting existing rows 95-117 down
# to 96-118), populating the new row with a fresh "Ají" price record for
# Vega Monumental Concepción. This matches the weekly refresh described in
# the commit message "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 95 (and everything below it) down by one to make room
# for the new record.
$ws.Rows("95:95").Insert(-4121)   # xlShiftDown = -4121

# Populate the freshly inserted row 95. Columns A, B, C, E, F, G, I and R
# are constant across this block of rows, so copy them from row 96 (the
# row that used to be row 95 before the insert). Use Value2 for reads,
# since the Value getter in this runtime does not return cell contents.
$ws.Range("A95").Value = $ws.Range("A96").Value2
$ws.Range("B95").Value = $ws.Range("B96").Value2
$ws.Range("C95").Value = $ws.Range("C96").Value2
$ws.Range("E95").Value = $ws.Range("E96").Value2
$ws.Range("F95").Value = $ws.Range("F96").Value2
$ws.Range("G95").Value = $ws.Range("G96").Value2
$ws.Range("I95").Value = $ws.Range("I96").Value2
$ws.Range("R95").Value = $ws.Range("R96").Value2

# New data specific to this record. D is formatted as a date (same number
# format as the rest of the "Fecha" column).
$ws.Range("D95").Value = 44754
$ws.Range("D95").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("H95").Value = "Inferno"
$ws.Range("J95").Value = 40
$ws.Range("K95").Value = 12000
$ws.Range("L95").Value = 13000
$ws.Range("M95").Value = 12500
$ws.Range("N95").Value = "$/caja 12 kilos"
$ws.Range("O95").Value = "Región de Arica y Parinacota"
$ws.Range("P95").Value = 1042
$ws.Range("Q95").Value = 12
